# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Writes the recalculated K values (and leaves std/mean unaffected since
# this sheet only stores the per-game values) for rows 2-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 8
    3  = 2
    4  = 4
    5  = 5
    6  = 8
    7  = 3
    8  = 2
    9  = 5
    10 = 4
    11 = 1
    12 = 7
    13 = 6
    14 = 7
    15 = 3
    16 = 3
    17 = 5
    18 = 5
    19 = 5
    20 = 0
    21 = 3
    22 = 3
    23 = 2
    24 = 6
    25 = 8
    26 = 0
    27 = 5
    28 = 4
    29 = 7
    30 = 3
    31 = 1
    32 = 6
    33 = 2
    34 = 2
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
